$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The source data set gained one more weekly reporting date (2021-09-09,
# serial 44448) for "Pera" at "Mercado Mayorista Lo Valledor de Santiago".
# Two new observations (Packham's Triumph, quality "Especial" and
# "Primera") are inserted as new rows 638-639, pushing the existing rows
# 638-657 down to 640-659.

# Insert two blank rows at 638:639 - Excel shifts rows 638-657 down to
# 640-659 and carries the row-638 formatting (date style) onto the new
# rows, matching the target dimension A1:T659.
$ws.Rows("638:639").Insert()

# --- New row 638 ---
$ws.Range("A638").Value = 6
$ws.Range("B638").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C638").Value = "Metropolitana"
$ws.Range("D638").Value = 44448
$ws.Range("E638").Value = 13
$ws.Range("F638").Value = "Fruta"
$ws.Range("G638").Value = 100104
$ws.Range("H638").Value = "Frutos de pepita"
$ws.Range("I638").Value = 100104005
$ws.Range("J638").Value = "Pera"
$ws.Range("K638").Value = "Packham's Triumph"
$ws.Range("L638").Value = "Especial"
$ws.Range("M638").Value = 18
$ws.Range("N638").Value = 190000
$ws.Range("O638").Value = 190000
$ws.Range("P638").Value = 190000
$ws.Range("Q638").Value = "$/bins (450 kilos)"
$ws.Range("R638").Value = "Región de O'Higgins"
$ws.Range("S638").Value = 422
$ws.Range("T638").Value = 450

# --- New row 639 ---
$ws.Range("A639").Value = 6
$ws.Range("B639").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C639").Value = "Metropolitana"
$ws.Range("D639").Value = 44448
$ws.Range("E639").Value = 13
$ws.Range("F639").Value = "Fruta"
$ws.Range("G639").Value = 100104
$ws.Range("H639").Value = "Frutos de pepita"
$ws.Range("I639").Value = 100104005
$ws.Range("J639").Value = "Pera"
$ws.Range("K639").Value = "Packham's Triumph"
$ws.Range("L639").Value = "Primera"
$ws.Range("M639").Value = 26
$ws.Range("N639").Value = 160000
$ws.Range("O639").Value = 170000
$ws.Range("P639").Value = 165000
$ws.Range("Q639").Value = "$/bins (450 kilos)"
$ws.Range("R639").Value = "Región de O'Higgins"
$ws.Range("S639").Value = 367
$ws.Range("T639").Value = 450
